$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing personal numbers
$ws.Range("A2").Value = 8
$ws.Range("A3").Value = 9

# Reset the explicit row height (30) back to the sheet's default row height
$ws.Rows("1:3").AutoFit()

# Add new row of data
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Pedri"
$ws.Range("C4").Value = "Gonzales"
$ws.Range("D4").Value = "Gràcia"
$ws.Range("E4").Value = "Barcelona"

# Copy the style from row 3 onto the new row 4
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to match the final state
$ws.Range("H11").Select()
